$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: new data row, following the same pattern as rows 3-15
# Column A uses the bold/bordered/centered style used by the rest of column A
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.031533375272873
$ws.Range("D16").Value = 0.9183273690739626
$ws.Range("E16").Value = 1.008475206288345
$ws.Range("F16").Value = 1.031533375272873
$ws.Range("G16").Value = 0.9473117055122916
$ws.Range("H16").Value = 1.038744308907424
$ws.Range("I16").Value = 1.013762454306349
$ws.Range("J16").Value = 0.9183273690739626
$ws.Range("K16").Value = 0.9634012876811539
$ws.Range("L16").Value = 0.9974673314770136
$ws.Range("M16").Value = 0.9930257365602078
